$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old summary rows (blank spacer row 100, SUM row 101) need to move down
# by one row to make room for a new card entry. Inserting a row right before
# the blank spacer row (old row 100) pushes everything after it down, and
# picks up the formatting of the row above (row 99 -> style index 3) the
# same way Excel does when a row is inserted mid-table.
$ws.Rows.Item(100).Insert()

# Fill in the new card row.
$ws.Range("A100").Value = "Palafin"
$ws.Range("B100").Value = "SV Black Star Promos"
$ws.Range("C100").Value = "SVP036"
$ws.Range("D100").Value = 4.2

# Update the SUM formula (now on row 102) to include the new row.
$ws.Range("D102").Formula = "=SUM(D2:D100)"

# Update the view so selection/scroll position matches the new extent.
$ws.Range("A70").Select()
$ws.Application.ActiveWindow.ScrollRow = 70
$ws.Range("D101").Select()
